$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the two added columns (System_Id / Device_Id)
$ws.Cells.Item(1, 6).Value = "System_Id"
$ws.Cells.Item(1, 7).Value = "Device_Id"

# Row 2 - NumericID 3200000 (temp_amb)
$ws.Cells.Item(2, 6).Value = "MET51"
$ws.Cells.Item(2, 7).Value = "B870"

# Row 3 - NumericID 3200008 (irrad)
$ws.Cells.Item(3, 6).Value = "MET51"
$ws.Cells.Item(3, 7).Value = "B870"

# Row 4 - NumericID 42180043 (Y1_meeting)
$ws.Cells.Item(4, 6).Value = "U25M2"
$ws.Cells.Item(4, 7).Value = "Y700"

# Row 5 - NumericID 42180045 (Y1_office)
$ws.Cells.Item(5, 6).Value = "U25M2"
$ws.Cells.Item(5, 7).Value = "Y701"

# Row 6 - NumericID 42180159 (blinds_height_F1)
$ws.Cells.Item(6, 6).Value = "U25F1"
$ws.Cells.Item(6, 7).Value = "M200"

# Row 7 - NumericID 42180160 (blinds_angle_F1)
$ws.Cells.Item(7, 6).Value = "U25F1"
$ws.Cells.Item(7, 7).Value = "M200"

# Row 8 - NumericID 42180170 (blinds_height_F2)
$ws.Cells.Item(8, 6).Value = "U25F2"
$ws.Cells.Item(8, 7).Value = "M200"

# Row 9 - NumericID 42180171 (blinds_angle_F2)
$ws.Cells.Item(9, 6).Value = "U25F2"
$ws.Cells.Item(9, 7).Value = "M200"

# Row 10 - NumericID 42180179 (blinds_height_F3)
$ws.Cells.Item(10, 6).Value = "U25F3"
$ws.Cells.Item(10, 7).Value = "M200"

# Row 11 - NumericID 42180187 (blinds_height_F4)
$ws.Cells.Item(11, 6).Value = "U25F4"
$ws.Cells.Item(11, 7).Value = "M200"

# Row 12 - NumericID 42180023 (temp_meeting)
$ws.Cells.Item(12, 6).Value = "U25M2"
$ws.Cells.Item(12, 7).Value = "B810"

# Row 13 - NumericID 42180034 (temp_office)
$ws.Cells.Item(13, 6).Value = "U25M2"
$ws.Cells.Item(13, 7).Value = "B811"

# Row 14 - NumericID 42180024 (setp_meeting)
$ws.Cells.Item(14, 6).Value = "U25M2"
$ws.Cells.Item(14, 7).Value = "B810"

# Row 15 - NumericID 42180035 (setp_office)
$ws.Cells.Item(15, 6).Value = "U25M2"
$ws.Cells.Item(15, 7).Value = "B811"

# Row 16 - NumericID 42160122 (heating_power)
$ws.Cells.Item(16, 6).Value = "U25M1"
$ws.Cells.Item(16, 7).Value = "P890"

# Row 17 - NumericID 42160203 (cooling_power)
$ws.Cells.Item(17, 6).Value = "U25N1"
$ws.Cells.Item(17, 7).Value = "P890"

# Row 18 - NumericID 42160278 (pv_active_power)
$ws.Cells.Item(18, 6).Value = "U25E3"
$ws.Cells.Item(18, 7).Value = "T100"

# Row 19 - NumericID 42160255 (total_active_power)
$ws.Cells.Item(19, 6).Value = "U25E1"
$ws.Cells.Item(19, 7).Value = "P001"

# Row 20 - NumericID 42180007 (praes_office)
$ws.Cells.Item(20, 6).Value = "U25R2"
$ws.Cells.Item(20, 7).Value = "B870"

# Row 21 - NumericID 42180000 (praes_meeting)
$ws.Cells.Item(21, 6).Value = "U25R1"
$ws.Cells.Item(21, 7).Value = "B870"

# Match the author's final selection position
$ws.Range("G16").Select()
